# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update the "K" column (G) values for each row with newly (re)calculated strike counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 0
    6  = 1
    7  = 1
    8  = 4
    10 = 1
    11 = 1
    12 = 0
    13 = 1
    14 = 0
    15 = 0
    16 = 3
    17 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
